$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cell types")

# Customer name correction (B5, next to "Musteri/Customer" label)
$ws.Range("B5").Value = "ege"

# Inspection place correction (B7, next to "Test Yeri/Inspection Place" label)
$ws.Range("B7").Value = "istanbul"

# Personal Bilgileri (personnel info) table, rows 34-37
# Row 34: Adi Soyadi / Name Surname -> fill in names
$ws.Range("B34").Value = "Ege"
$ws.Range("D34").Value = "Barış"
$ws.Range("E34").Value = "Serra"

# Row 35: Seviye / Level -> fill in levels, using copy/paste of existing
# text-typed numeric cells so the values land as shared-string text
# (matching the rest of the sheet) instead of numeric cells.
$ws.Range("E28").Copy()
$ws.Range("B35").PasteSpecial(-4163)
$ws.Range("D28").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E28").Copy()
$ws.Range("E35").PasteSpecial(-4163)

$excel.CutCopyMode = 0
